$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36
$ws.Range("A33:E33").Copy()
$ws.Range("A36:E36").PasteSpecial(-4122)
$ws.Range("A36").Value = "SCRIPT/T01P01A/um2101.ssb"
$ws.Range("B36").Value = 483
$ws.Range("C36").Value = " I\'ve been telling everyone the\ntruth, by golly!"
$ws.Range("D36").Value = " Ей-богу, я всем рассказываю\nправду!"
$ws.Range("E36").Value = " Åê-áïãô, ÿ âòåí ñàòòëàèúâàý\nðñàâäô!"
$ws.Rows.Item(36).RowHeight = 43.2

# Row 37
$ws.Range("A33:E33").Copy()
$ws.Range("A37:E37").PasteSpecial(-4122)
$ws.Range("A37").Clear()
$ws.Range("B37").Value = 486
$ws.Range("C37").Value = " Folks find it a hard tale to\nswallow. But I\'m winning them over with my\nhonesty."
$ws.Range("D37").Value = " Народу трудно в неё поверить.\nНо я убеждаю их своей прямотой."
$ws.Range("E37").Value = " Îàñïäô óñôäîï â îåæ ðïâåñéóû.\nÎï ÿ ôáåçäàý éö òâïåê ðñÿíïóïê."
$ws.Rows.Item(37).RowHeight = 31.8

# Row 38
$ws.Range("A34:E34").Copy()
$ws.Range("A38:E38").PasteSpecial(-4122)
$ws.Range("B38").Value = 489
$ws.Range("C38").Value = " I\'m surely getting through to\nfolks, yup yup!"
$ws.Range("D38").Value = " Я достучусь до каждого, да-да!"
$ws.Range("E38").Value = " Ÿ äïòóôœôòû äï ëàçäïãï, äà-äà!"
$ws.Rows.Item(38).RowHeight = 21.6

# Row 39
$ws.Range("A31:E31").Copy()
$ws.Range("A39:E39").PasteSpecial(-4122)
$ws.Range("A39").Value = "SCRIPT/T01P02A/um2207.ssb"
$ws.Range("B39").Value = 464
$ws.Range("C39").Value = " I\'m getting kitted out to prepare\nto go to [CS:P]Brine Cave[CR]... Huff-huff…"
$ws.Range("D39").Value = " Я снаряжаюсь, чтобы подготовиться\nк походу в [CS:P]Пещеру у Моря[CR]... Ух-ух..."
$ws.Range("E39").Value = " Ÿ òîàñÿçàýòû, œóïáú ðïäãïóïâéóûòÿ\në ðïöïäô â [CS:P]Ðåþåñô ô Íïñÿ[CR]... Ôö-ôö..."
$ws.Rows.Item(39).RowHeight = 43.2

# Row 40
$ws.Range("A33:E33").Copy()
$ws.Range("A40:E40").PasteSpecial(-4122)
$ws.Range("A40").Value = "SCRIPT/D25P11A/um2303.ssb"
$ws.Range("B40").Value = 445
$ws.Range("C40").Value = " Golly, this surely isn\'t easy.[K]\nBut I\'m not giving up, no sirree!"
$ws.Range("D40").Value = " Ей-богу, это совсем непросто.[K]\nНо я не сдамся, нет, сэррр!"
$ws.Range("E40").Value = " Åê-áïãô, üóï òïâòåí îåðñïòóï.[K]\nÎï ÿ îå òäàíòÿ, îåó, òüñññ!"
$ws.Rows.Item(40).RowHeight = 43.2

# Row 41
$ws.Range("A34:E34").Copy()
$ws.Range("A41:E41").PasteSpecial(-4122)
$ws.Range("A41").Value = "SCRIPT/D25P11A/um2306.ssb"
$ws.Rows.Item(41).RowHeight = 43.2

# Row 42
$ws.Range("A33:E33").Copy()
$ws.Range("A42:E42").PasteSpecial(-4122)
$ws.Range("A42").Value = "SCRIPT/G01P03A/um2401.ssb"
$ws.Range("B42").Value = 420
$ws.Range("C42").Value = " Oof...[K] I can\'t do anything for\nyou all but give encouragement. It frustrates\nme something awful!"
$ws.Range("D42").Value = " Ууф...[K] Я ничего не могу для вас\nсделать, только поддержать. Это меня дико\nрасстраивает!"
$ws.Range("E42").Value = " Ôôõ...[K] Ÿ îéœåãï îå íïãô äìÿ âàò\nòäåìàóû, óïìûëï ðïääåñçàóû. Üóï íåîÿ äéëï\nñàòòóñàéâàåó!"
$ws.Rows.Item(42).RowHeight = 43.2

# Row 43
$ws.Range("A33:E33").Copy()
$ws.Range("A43:E43").PasteSpecial(-4122)
$ws.Range("A43").Clear()
$ws.Range("B43").Value = 423
$ws.Range("C43").Value = " But...[K]you all get it done now!"
$ws.Range("D43").Value = " Но...[K] Вы справитесь!"
$ws.Range("E43").Value = " Îï...[K] Âú òðñàâéóåòû!"

# Row 44
$ws.Range("A33:E33").Copy()
$ws.Range("A44:E44").PasteSpecial(-4122)
$ws.Range("A44").Clear()
$ws.Range("B44").Value = 426
$ws.Range("C44").Value = " You all don\'t fail us now! Stop\nour planet from becoming paralyzed! You hear?"
$ws.Range("D44").Value = " Не подведите нас! Остановите\nпланетарный паралич! Слышите?"
$ws.Range("E44").Value = " Îå ðïäâåäéóå îàò! Ïòóàîïâéóå\nðìàîåóàñîúê ðàñàìéœ! Òìúšéóå?"
$ws.Rows.Item(44).RowHeight = 31.8

$excel.Application.CutCopyMode = $false
$ws.Range("E44").Select()
